# Update specific numeric values in the result_data_KNN sheet.
# These correspond to algorithm result changes ("Update Name of Algo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.782
$ws.Range("C9").Value = -11.445
$ws.Range("E11").Value = 12.914
$ws.Range("C18").Value = -12.397
$ws.Range("C20").Value = -12.581
